# Updated analyses with extended species:
# refresh the lnL / omega / LRT table on Sheet1 with the re-run values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -13005.275962
$ws.Range("C2").Value = 0.21745

$ws.Range("B3").Value = -13003.538552
$ws.Range("C3").Value = 0.21965000000000001
$ws.Range("D3").Value = 0.077619999999999995

$ws.Range("B4").Value = -12871.89337
$ws.Range("C4").Value = 0.092050000000000007
$ws.Range("D4").Value = 0.41310000000000002

$ws.Range("B5").Value = -12867.384028
$ws.Range("C5").Value = 0.11772000000000001
$ws.Range("D5").Value = 0.41298000000000001
$ws.Range("E5").Value = 0.075029999999999999

$ws.Range("C9").Activate()
